$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F header: from_pvt1 -> from_heat_pump1 (heat pump "from" source added; CHP column removed
# by shifting CHP1(col E) out and pvt1 into col E, heat_pump1 into col F)
$ws.Range("E1").Value = "P_from_pvt1"
$ws.Range("F1").Value = "P_from_heat_pump1"

# Row 2 (to demand1)
$ws.Range("E2").Value = "P_pvt1_demand1"
$ws.Range("F2").Value = 0

# Row 3 (to net1)
$ws.Range("E3").Value = "P_pvt1_net1"
$ws.Range("F3").Value = 0

# Row 4 (to bat1)
$ws.Range("E4").Value = "P_pvt1_bat1"
$ws.Range("F4").Value = 0

# Row 5 (to charging_station1) - bat1 now feeds charging_station1
$ws.Range("D5").Value = "P_bat1_charging_station1"
$ws.Range("E5").Value = "P_pvt1_charging_station1"
$ws.Range("F5").Value = 0

# Row 6: charging_station2 -> heat_pump1
$ws.Range("A6").Value = "P_to_heat_pump1"
$ws.Range("B6").Value = "P_net1_heat_pump1"
$ws.Range("C6").Value = "P_pv1_heat_pump1"
$ws.Range("D6").Value = "P_bat1_heat_pump1"
$ws.Range("E6").Value = "P_pvt1_heat_pump1"
$ws.Range("F6").Value = 0
